# "Gift card page added"
#
# 1. Insert a new "GiftCard" worksheet right before "Test" (Worksheets.Add()
#    inserts before the active sheet, and "Test" is the active/tabSelected
#    sheet in the source workbook, so this lands it in the right spot and
#    picks up sheetId 10 while "Test" keeps its original sheetId 9).
# 2. Seed GiftCard's formatting by copying the old Test header/format block
#    (A1:G3) into it, then overwrite the values for the 7 gift-card rows.
# 3. Trim the "Test" sheet back down to the small 7-column / 2-row table it
#    keeps after the edit.
# 4. Minor cosmetic view-state changes (workbook window, Customer selection).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new GiftCard sheet (goes in before "Test" automatically).
# ---------------------------------------------------------------------
$testSheet = $wb.Worksheets.Item("Test")
$giftCard = $wb.Worksheets.Add()
$giftCard.Name = "GiftCard"

# Re-fetch a live reference to Test (stale after Worksheets.Add()).
$testSheet = $wb.Worksheets.Item("Test")

# Copy the old Test sheet's A1:G3 block (values + number formats/fills/
# borders) into GiftCard so we inherit the right style indices, then we'll
# overwrite every cell with the actual GiftCard content/styling below.
$testSheet.Range("A1:G3").Copy($giftCard.Range("A1"))

# Re-fetch again after the copy.
$giftCard = $wb.Worksheets.Item("GiftCard")
$testSheet = $wb.Worksheets.Item("Test")

# ---------------------------------------------------------------------
# 2) Fix up column D (expiry date-as-text) and column G (customer name)
#    styles by pulling in the exact existing style ("s=20" / "s=2") from
#    elsewhere in the workbook, then stamp rows 4-7 with the same look as
#    row 2 (s19/s23/s2/s20 combo) via a formats-only paste so we reuse
#    existing style records instead of minting new ones.
# ---------------------------------------------------------------------
$testSheet.Range("AF2").Copy()
$giftCard.Range("D2:D3").PasteSpecial(-4122)

$customerSheet = $wb.Worksheets.Item("Customer")
$customerSheet.Range("A2").Copy()
$giftCard.Range("G2:G3").PasteSpecial(-4122)

$giftCard.Range("A2:G2").Copy()
$giftCard.Range("A4:G4").PasteSpecial(-4122)
$giftCard.Range("A5:G5").PasteSpecial(-4122)
$giftCard.Range("A6:G6").PasteSpecial(-4122)
$giftCard.Range("A7:G7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Header row (left to right) + data columns (top to bottom) - this
#    column-major fill order for the data rows matches the shared-string
#    allocation order seen in the authored workbook ("Promotional" before
#    the "5/31/2025" / "6/30/2025" expiry strings).
# ---------------------------------------------------------------------
$giftCard.Range("A1").Value = "Type"
$giftCard.Range("B1").Value = "Number"
$giftCard.Range("C1").Value = "Reason"
$giftCard.Range("D1").Value = "Expiry"
$giftCard.Range("E1").Value = "Status"
$giftCard.Range("F1").Value = "Balance"
$giftCard.Range("G1").Value = "Customer"

$giftCard.Range("A2").Value = "Normal"
$giftCard.Range("A3").Value = "Normal"
$giftCard.Range("A4").Value = "Normal"
$giftCard.Range("A5").Value = "Promotional"
$giftCard.Range("A6").Value = "Promotional"
$giftCard.Range("A7").Value = "Promotional"

$giftCard.Range("B2").Value = 1001
$giftCard.Range("B3").Value = 1002
$giftCard.Range("B4").Value = 1003
$giftCard.Range("B5").Value = 1004
$giftCard.Range("B6").Value = 1005
$giftCard.Range("B7").Value = 1006

$giftCard.Range("C2").Value = "Store Credit"
$giftCard.Range("C3").Value = "Open House Promo"
$giftCard.Range("C4").Value = "Store Credit"
$giftCard.Range("C5").Value = "Store Credit"
$giftCard.Range("C6").Value = "Store Credit"
$giftCard.Range("C7").Value = "Store Credit"

$giftCard.Range("D2").Value = "5/31/2025"
$giftCard.Range("D3").Value = "6/30/2025"
$giftCard.Range("D4").Value = "5/31/2025"
$giftCard.Range("D5").Value = "6/30/2025"
$giftCard.Range("D6").Value = "5/31/2025"
$giftCard.Range("D7").Value = "6/30/2025"

$giftCard.Range("E2").Value = "Active"
$giftCard.Range("E3").Value = "Active"
$giftCard.Range("E4").Value = "Inactive"
$giftCard.Range("E5").Value = "Active"
$giftCard.Range("E6").Value = "Active"
$giftCard.Range("E7").Value = "Suspended"

$giftCard.Range("F2").Value = 500
$giftCard.Range("F3").Value = 500
$giftCard.Range("F4").Value = 0
$giftCard.Range("F5").Value = 1000
$giftCard.Range("F6").Value = 1000
$giftCard.Range("F7").Value = 0

$giftCard.Range("G2").Value = ""
$giftCard.Range("G3").Value = "Maria Lopez"
$giftCard.Range("G4").Value = "Emily Johnson"

$giftCard.Rows.Item(1).RowHeight = 16.8

# Column widths approximating the authored sheet.
$giftCard.Columns.Item(1).ColumnWidth = 15.83
$giftCard.Columns.Item(2).ColumnWidth = 15.94
$giftCard.Columns.Item(3).ColumnWidth = 21.83
$giftCard.Columns.Item(4).ColumnWidth = 13.39
$giftCard.Columns.Item(5).ColumnWidth = 15.5
$giftCard.Columns.Item(6).ColumnWidth = 17.83
$giftCard.Columns.Item(7).ColumnWidth = 28.39

$giftCard.Range("L14").Select()

# ---------------------------------------------------------------------
# 3) Trim "Test" down to the small 7-column / 2-row table.
# ---------------------------------------------------------------------
$testSheet.Columns("H:AL").Delete()

$testSheet.Range("A1").Value = "Type"
$testSheet.Range("B1").Value = "Number"
$testSheet.Range("C1").Value = "Reason"
$testSheet.Range("D1").Value = "Expiry"
$testSheet.Range("E1").Value = "Status"
$testSheet.Range("F1").Value = "Balance"
$testSheet.Range("G1").Value = "Customer"

$testSheet.Range("A2").Value = "Normal"
$testSheet.Range("B2").Value = 2002
$testSheet.Range("C2").Value = "Open House Promo"
$testSheet.Range("D2").Value = "6/30/2025"
$testSheet.Range("E2").Value = "Active"
$testSheet.Range("F2").Value = 500
$testSheet.Range("G2").Value = "Maria Lopez"

$testSheet.Range("A3:G3").Clear()

$testSheet.Columns.Item(1).ColumnWidth = 17.17
$testSheet.Columns.Item(7).ColumnWidth = 10.39

$testSheet.Range("L14").Select()

# ---------------------------------------------------------------------
# 4) Cosmetic view-state tweaks.
# ---------------------------------------------------------------------
$customerSheet = $wb.Worksheets.Item("Customer")
$customerSheet.Range("C2:C7").Select()

$excel.Windows.Item(1).WindowState = -4143
